$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-9 from 45174 to 45175
$ws.Range("C2:C9").Value = 45175
